$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.774.22'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.493.73'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '532.69'
$ws.Range('E5').Value = '  +4.10%  '
$ws.Range('D6').Value = '133.63'
$ws.Range('E6').Value = '  +2.93%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').Value = '  +3.63%  '
$ws.Range('D9').Value = '2.502.35'
$ws.Range('E10').Value = '  +3.13%  '
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '2.935.47'
$ws.Range('E14').Value = '  +1.74%  '
$ws.Range('D15').Value = '58.616.93'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '22.24'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = '2.495.22'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +2.98%  '
$ws.Range('D21').Value = '320.42'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').Value = '6.20'
$ws.Range('E22').Value = '  +4.88%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '66.09'
$ws.Range('E24').Value = '  +4.96%  '
$ws.Range('D25').Value = '0.406'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('D26').Value = '0.991'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').Value = '7.44'
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('D29').Value = '173.04'
$ws.Range('E29').Value = '  +2.60%  '
$ws.Range('D30').Value = '0.0₃0755'
$ws.Range('E30').Value = '  +3.40%  '
$ws.Range('E31').Value = '  +3.78%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = '6.28'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '18.07'
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').Value = '0.835'
$ws.Range('E39').Value = '  +9.58%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.50'
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('D41').Value = '36.34'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('E42').Value = '  +2.79%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '274.46'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '5.05'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '131.09'
$ws.Range('E45').Value = '  +8.90%  '
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('E48').Value = '  +4.47%  '
$ws.Range('E49').Value = '  +3.30%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '16.77'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.752.25'
$ws.Range('E51').Value = '  +2.59%  '
